# Sync attendance_reports: normalize "Recorded By" (column G) ordering.
# For every data row, the comma-separated list of recorders in column G
# has its token order reversed (e.g. "a, b" -> "b, a").
# Single-token values are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -gt 1) {
            $reversed = $parts[($parts.Length - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
